# UC014 - Minha Conta Bancária : version bump + wording/punctuation fixes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Version: 1.0 -> 1.2.5
$ws.Range("D2").Value = "1.2.5"

# Precondition text: fix "usuario" -> "usuário", add trailing period
$ws.Range("B8").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B17").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B27").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B37").Value = "O usuário devidamente autenticado e na tela inicial do sistema."

# Step 1 text: add trailing period
$ws.Range("B10").Value = "Chefe Acessa a funcionalidade Minha Conta Bancária (menu)."
$ws.Range("B19").Value = "Chefe Acessa a funcionalidade Minha Conta Bancária (menu)."
$ws.Range("B29").Value = "Chefe Acessa a funcionalidade Minha Conta Bancária (menu)."
$ws.Range("B39").Value = "Chefe Acessa a funcionalidade Minha Conta Bancária (menu)."

# Expected result text: "bancários" -> "bancária"
$ws.Range("D11").Value = "SYSTEM Exibe mensagens informativas (MSG403 - Informativos sobre a atualização de conta bancária (dados bancários)) para o usuário sobre a manutenção de informações bancárias."

# Expected result text: add trailing period
$ws.Range("D20").Value = "SYSTEM Apresenta os campos (banco/agência/conta corrente) alterados."
$ws.Range("D30").Value = "SYSTEM Apresenta os campos (banco/agência/conta corrente) alterados."
$ws.Range("D40").Value = "SYSTEM Apresenta os campos (banco/agência/conta corrente) alterados."
